# Actualización de la planeación del sistema
# Inserta tres nuevas actividades (Creación de C.U., Creación de prototipos,
# Integración de CSS a sistema) antes de "Vista de horarios de materias (maestros)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows before the current row 9, pushing existing rows 9-16 down to 12-19.
$ws.Rows("9:11").Insert()

# Row 9: Creación de C.U.
$ws.Range("A9").Value = "Creación de C.U."
$ws.Range("B9").Value = "Beatriz"
$ws.Range("C9").Value = "Domingo 3 de mayo"
$ws.Range("D9").Value = "Domingo 3 de mayo"
$ws.Range("E9").Value = "Martes 5 de mayo"

# Row 10: Creación de prototipos
$ws.Range("A10").Value = "Creación de prototipos"
$ws.Range("B10").Value = "Julian"
$ws.Range("C10").Value = "Domingo 3 de mayo"
$ws.Range("D10").Value = "Domingo 3 de mayo"
$ws.Range("E10").Value = "Martes 5 de mayo"

# Row 11: Integración de CSS a sistema
$ws.Range("A11").Value = "Integración de CSS a sistema"
$ws.Range("B11").Value = "Jesús"
$ws.Range("C11").Value = "Domingo 3 de mayo"
$ws.Range("D11").Value = "Domingo 3 de mayo"
$ws.Range("E11").Value = "Martes 5 de mayo"

# Row 12 (formerly row 9, "Vista de horarios de materias (maestros)") gets an
# updated delivery date.
$ws.Range("E12").Value = "Martes 12 de mayo"

# Update the selected cell to match the saved view state.
$ws.Range("A8").Select()
